# Add two new columns, I0 (col I) and IF (col J), to the right of the
# existing IP column (col H), mirroring the header style already used
# by the other header cells in row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - copy the style from the existing header cell H1 so the
# new headers match (bold, centered, bordered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-21 for columns I (I0) and J (IF)
$values = @{
    2  = @(1, 4)
    3  = @(1, 5)
    4  = @(1, 7)
    5  = @(1, 4)
    6  = @(7, 8)
    7  = @(7, 8)
    8  = @(5, 8)
    9  = @(1, 5)
    10 = @(1, 4)
    11 = @(7, 7)
    12 = @(7, 8)
    13 = @(1, 6)
    14 = @(1, 5)
    15 = @(2, 6)
    16 = @(7, 8)
    17 = @(1, 4)
    18 = @(1, 4)
    19 = @(4, 5)
    20 = @(6, 7)
    21 = @(3, 4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
